$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.01155686378479
$ws.Range("B1").Value = 2.129428625106812
$ws.Range("C1").Value = 5.758969306945801
$ws.Range("D1").Value = 0.9789129495620728
$ws.Range("E1").Value = 1.060928702354431
